$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.086.43'
$ws.Range('E2').Value = '  -0.57%  '
$ws.Range('D3').Value = '2.511.85'
$ws.Range('E3').Value = '  -0.35%  '
$ws.Range('D4').Value = "'0.999"
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').Value = "'537.07"
$ws.Range('E5').Value = '  -0.92%  '
$ws.Range('D6').Value = "'136.95"
$ws.Range('E6').Value = '  -2.14%  '
$ws.Range('E7').Value = '  +0.14%  '
$ws.Range('E8').Value = '  +0.62%  '
$ws.Range('D9').Value = '2.525.96'
$ws.Range('E9').Value = '  -0.01%  '
$ws.Range('E10').Value = '  +0.06%  '
$ws.Range('E11').Value = '  -2.54%  '
$ws.Range('D12').Value = "'5.30"
$ws.Range('E12').Value = '  -2.27%  '
$ws.Range('E13').Value = '  -0.58%  '
$ws.Range('D14').Value = '2.973.84'
$ws.Range('E14').Value = '  +0.00%  '
$ws.Range('D15').Value = "'23.11"
$ws.Range('E15').Value = '  -0.91%  '
$ws.Range('D16').Value = '59.138.26'
$ws.Range('E16').Value = '  -0.33%  '
$ws.Range('E17').Value = '  -1.58%  '
$ws.Range('D18').Value = '2.525.02'
$ws.Range('E18').Value = '  +0.31%  '
$ws.Range('E19').Value = '  +0.46%  '
$ws.Range('E20').Value = '  +0.01%  '
$ws.Range('D21').Value = "'323.96"
$ws.Range('E22').Value = '  +0.09%  '
$ws.Range('E23').Value = '  +1.19%  '
$ws.Range('D24').Value = "'65.45"
$ws.Range('E24').Value = '  +3.38%  '
$ws.Range('E25').Value = '  -0.09%  '
$ws.Range('E26').Value = '  -1.70%  '
$ws.Range('E27').Value = '  -0.12%  '
$ws.Range('E28').Value = '  -2.46%  '
$ws.Range('E29').Value = '  -1.49%  '
$ws.Range('E30').Value = '  -0.43%  '
$ws.Range('D31').Value = "'172.05"
$ws.Range('E31').Value = '  +3.98%  '
$ws.Range('D32').Value = "'1.77"
$ws.Range('E32').Value = '  -1.66%  '
$ws.Range('E33').Value = '  +5.50%  '
$ws.Range('E35').Value = '  +0.70%  '
$ws.Range('D36').Value = "'18.41"
$ws.Range('E36').Value = '  -0.90%  '
$ws.Range('E37').Value = '  -1.13%  '
$ws.Range('E38').Value = '  -3.08%  '
$ws.Range('D39').Value = "'36.74"
$ws.Range('E39').Value = '  -0.68%  '
$ws.Range('E40').Value = '  +0.44%  '
$ws.Range('E41').Value = '  -2.03%  '
$ws.Range('D42').Value = "'285.48"
$ws.Range('E42').Value = '  +1.45%  '
$ws.Range('D43').Value = "'5.12"
$ws.Range('E43').Value = '  -2.15%  '
$ws.Range('B44').Value = 'Mantle'
$ws.Range('C44').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D44').Value = "'0.612"
$ws.Range('E44').Value = '  +1.77%  '
$ws.Range('B45').Value = 'FirstDigitalUSD'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D45').Value = "'0.997"
$ws.Range('E45').Value = '  -0.06%  '
$ws.Range('D46').Value = "'131.62"
$ws.Range('E46').Value = '  +4.18%  '
$ws.Range('E47').Value = '  -0.02%  '
$ws.Range('E48').Value = '  -1.65%  '
$ws.Range('D49').Value = "'0.0508"
$ws.Range('E49').Value = '  -0.98%  '
$ws.Range('E50').Value = '  -1.34%  '
$ws.Range('D51').Value = "'17.46"
$ws.Range('E51').Value = '  -2.34%  '
